$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Heirloom Tomato): update description text
$ws.Range("C3").Value = "Not your standard tomato"

# Row 3 (Heirloom Tomato continued): make the Salt column (H3) explicit value 0
$ws.Range("H3").Value = 0

# Row 4 (Spaghetti): add a description
$ws.Range("C4").Value = "Pairs with thick, tomato based sauces."

# Row 5 (Mom's Spaghetti): update description text
$ws.Range("C5").Value = "Always Ready, Mom's Spaghetti"

# Update the active selection to match the edit location
$ws.Range("H3").Select()
